# Auto-generated Excel COM-interop script applying numeric updates
# to the Leve profit-tracking sheets, per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 1502.9642
$ws.Range("I15").Value = 1502.9642
$ws.Range("K15").Value = 4508.892599999999
$ws.Range("M15").Value = -4339.892599999999
# ALC row 53
$ws.Range("H53").Value = 621.1818
$ws.Range("I53").Value = 529.625
$ws.Range("K53").Value = 529.625
$ws.Range("M53").Value = 107.375
# ALC row 76
$ws.Range("H76").Value = 2099.6667
$ws.Range("I76").Value = 1650
$ws.Range("J76").Value = 2999
$ws.Range("K76").Value = 1650
$ws.Range("L76").Value = 2999
$ws.Range("M76").Value = -1335
$ws.Range("N76").Value = -3629
# ALC row 79
$ws.Range("H79").Value = 2099.6667
$ws.Range("I79").Value = 1650
$ws.Range("J79").Value = 2999
$ws.Range("K79").Value = 1650
$ws.Range("L79").Value = 2999
$ws.Range("M79").Value = -558
$ws.Range("N79").Value = -5183
# ALC row 112
$ws.Range("H112").Value = 1005.41174
$ws.Range("J112").Value = 1084
$ws.Range("L112").Value = 3252
$ws.Range("N112").Value = -5468
# ALC row 123
$ws.Range("H123").Value = 180000
$ws.Range("J123").Value = 180000
$ws.Range("L123").Value = 180000
$ws.Range("N123").Value = -189800
# ALC row 137
$ws.Range("H137").Value = 1933.35
$ws.Range("I137").Value = 1333.5714
$ws.Range("J137").Value = 3332.8333
$ws.Range("K137").Value = 4000.7142
$ws.Range("L137").Value = 9998.499899999999
$ws.Range("M137").Value = -1450.7142
$ws.Range("N137").Value = -15098.4999
# ALC row 138
$ws.Range("H138").Value = 7519.0356
$ws.Range("I138").Value = 3831.3333
$ws.Range("J138").Value = 7961.56
$ws.Range("K138").Value = 11493.9999
$ws.Range("L138").Value = 23884.68
$ws.Range("M138").Value = -6353.999899999999
$ws.Range("N138").Value = -34164.68

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 3382.3447
$ws.Range("I32").Value = 2788.8928
$ws.Range("K32").Value = 2788.8928
$ws.Range("M32").Value = -2501.8928
# ARM row 45
$ws.Range("H45").Value = 2773.5
$ws.Range("I45").Value = 2365
$ws.Range("J45").Value = 3999
$ws.Range("K45").Value = 2365
$ws.Range("L45").Value = 3999
$ws.Range("M45").Value = -1988
$ws.Range("N45").Value = -4753
# ARM row 63
$ws.Range("H63").Value = 3755.4
$ws.Range("I63").Value = 3772
$ws.Range("J63").Value = 3689
$ws.Range("K63").Value = 3772
$ws.Range("L63").Value = 3689
$ws.Range("M63").Value = -3086
$ws.Range("N63").Value = -5061
# ARM row 66
$ws.Range("H66").Value = 3755.4
$ws.Range("I66").Value = 3772
$ws.Range("J66").Value = 3689
$ws.Range("K66").Value = 18860
$ws.Range("L66").Value = 18445
$ws.Range("M66").Value = -15428
$ws.Range("N66").Value = -25309
# ARM row 74
$ws.Range("H74").Value = 1922.5
$ws.Range("I74").Value = 1107
$ws.Range("K74").Value = 1107
$ws.Range("M74").Value = -233
# ARM row 77
$ws.Range("H77").Value = 1922.5
$ws.Range("I77").Value = 1107
$ws.Range("K77").Value = 5535
$ws.Range("M77").Value = -1167
# ARM row 97
$ws.Range("H97").Value = 778.75
$ws.Range("I97").Value = 854.2857
$ws.Range("J97").Value = 250
$ws.Range("K97").Value = 854.2857
$ws.Range("L97").Value = 250
$ws.Range("M97").Value = -358.2857
$ws.Range("N97").Value = -1242
# ARM row 110
$ws.Range("H110").Value = 3656.5
$ws.Range("I110").Value = 3775.6667
$ws.Range("J110").Value = 3299
$ws.Range("K110").Value = 3775.6667
$ws.Range("L110").Value = 3299
$ws.Range("M110").Value = -1730.6667
$ws.Range("N110").Value = -7389

$ws = $wb.Worksheets.Item("BSM")
# BSM row 105
$ws.Range("H105").Value = 2666
$ws.Range("I105").Value = 2999.5
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 2999.5
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -1252.5
$ws.Range("N105").Value = -5493
# BSM row 134
$ws.Range("H134").Value = 4211.6665
$ws.Range("I134").Value = 4211.6665
$ws.Range("K134").Value = 12634.9995
$ws.Range("M134").Value = -10099.9995

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 4005.5405
$ws.Range("I31").Value = 1553
$ws.Range("J31").Value = 9115
$ws.Range("K31").Value = 1553
$ws.Range("L31").Value = 9115
$ws.Range("M31").Value = -1258
$ws.Range("N31").Value = -9705
# CRP row 34
$ws.Range("H34").Value = 4005.5405
$ws.Range("I34").Value = 1553
$ws.Range("J34").Value = 9115
$ws.Range("K34").Value = 1553
$ws.Range("L34").Value = 9115
$ws.Range("M34").Value = -1351
$ws.Range("N34").Value = -9519
# CRP row 94
$ws.Range("H94").Value = 2012.1666
$ws.Range("I94").Value = 1816.6
$ws.Range("J94").Value = 2990
$ws.Range("K94").Value = 1816.6
$ws.Range("L94").Value = 2990
$ws.Range("M94").Value = -1365.6
$ws.Range("N94").Value = -3892
# CRP row 107
$ws.Range("H107").Value = 1149.5
$ws.Range("I107").Value = 763.8570999999999
$ws.Range("K107").Value = 763.8570999999999
$ws.Range("M107").Value = 1156.1429
# CRP row 132
$ws.Range("H132").Value = 5078.9
$ws.Range("I132").Value = 4448.5
$ws.Range("K132").Value = 13345.5
$ws.Range("M132").Value = -10815.5
# CRP row 133
$ws.Range("H133").Value = 124700
$ws.Range("J133").Value = 124700
$ws.Range("L133").Value = 124700
$ws.Range("N133").Value = -129760

$ws = $wb.Worksheets.Item("CUL")
# CUL row 11
$ws.Range("H11").Value = 5000
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 15000
$ws.Range("N11").Value = -15280
# CUL row 34
$ws.Range("H34").Value = 3131.6667
$ws.Range("J34").Value = 4830
$ws.Range("L34").Value = 14490
$ws.Range("N34").Value = -14658
# CUL row 128
$ws.Range("H128").Value = 250000
$ws.Range("I128").Value = 250000
$ws.Range("K128").Value = 750000
$ws.Range("M128").Value = -745020

$ws = $wb.Worksheets.Item("GSM")
# GSM row 80
$ws.Range("H80").Value = 5518.8335
$ws.Range("I80").Value = 5398.5557
$ws.Range("K80").Value = 5398.5557
$ws.Range("M80").Value = -4400.5557
# GSM row 83
$ws.Range("H83").Value = 5518.8335
$ws.Range("I83").Value = 5398.5557
$ws.Range("K83").Value = 26992.7785
$ws.Range("M83").Value = -22000.7785
# GSM row 122
$ws.Range("H122").Value = 999.6667
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897
# GSM row 126
$ws.Range("H126").Value = 2117.75
$ws.Range("I126").Value = 2117.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6353.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3883.25
$ws.Range("N126").Value = $null
# GSM row 132
$ws.Range("H132").Value = 5466.8
$ws.Range("I132").Value = 5238.857
$ws.Range("J132").Value = 5998.6665
$ws.Range("K132").Value = 15716.571
$ws.Range("L132").Value = 17995.9995
$ws.Range("M132").Value = -13186.571
$ws.Range("N132").Value = -23055.9995

$ws = $wb.Worksheets.Item("LTW")
# LTW row 55
$ws.Range("H55").Value = 459.75
$ws.Range("I55").Value = 384
$ws.Range("K55").Value = 384
$ws.Range("M55").Value = -211
# LTW row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = $null
# LTW row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = $null
# LTW row 82
$ws.Range("H82").Value = 1997.5
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 2000
$ws.Range("M82").Value = -1639
# LTW row 85
$ws.Range("H85").Value = 1997.5
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 2000
$ws.Range("M85").Value = -752

$ws = $wb.Worksheets.Item("WVR")
# WVR row 107
$ws.Range("H107").Value = 389.8
$ws.Range("I107").Value = 362.75
$ws.Range("K107").Value = 1088.25
$ws.Range("M107").Value = 831.75
# WVR row 132
$ws.Range("H132").Value = 2430.9412
$ws.Range("I132").Value = 2148.7693
$ws.Range("K132").Value = 6446.3079
$ws.Range("M132").Value = -3916.3079
